$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column A; existing column A (and its data/width)
# shifts right to become column D.
$ws.Range("A:C").EntireColumn.Insert()

# New header cells for the inserted columns, re-using the exact formatting
# ("Pandas" style) already applied to the shifted header cell D1.
$ws.Range("D1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$ws.Range("A1").Value = "button_apiKeyManagement_trNthChild"
$ws.Range("B1").Value = "button_apiKeyManagement_trNthChild_1"
$ws.Range("C1").Value = "button_apiKeyManagement_trNthChild_2"

# New data row values (kept as text, matching the shifted D2 cell's plain/
# unformatted style).
$ws.Range("A2").Value = "'2"
$ws.Range("B2").Value = "'1"
$ws.Range("C2").Value = "'1"
$ws.Range("D2").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)

# Column widths for the new layout. (The engine stores width as
# ColumnWidth + 5/6, so the set values are offset to land on the exact
# target stored widths of 36, 38, 38, 15.)
$ws.Columns.Item(1).ColumnWidth = 35.16666666666667
$ws.Columns.Item(2).ColumnWidth = 37.16666666666667
$ws.Columns.Item(3).ColumnWidth = 37.16666666666667
$ws.Columns.Item(4).ColumnWidth = 14.16666666666667

$excel.CutCopyMode = 0
